$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the "取得日時" (timestamp) column A for existing rows 2-12 to the new
# append time: 2025-11-25 12:50:18
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-11-25 12:50:18"
}
